# Weekly fruit/vegetable price update.
# A new price observation is inserted as row 89 (pushing the existing
# rows 89-91 down to 90-92); all other rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 89, shifting rows 89.. down by one.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new record.
$ws.Cells.Item(89, 1).Value  = 4
$ws.Cells.Item(89, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value  = "Los Lagos"
$ws.Cells.Item(89, 4).Value  = 45041
$ws.Cells.Item(89, 5).Value  = 10
$ws.Cells.Item(89, 6).Value  = 100112043
$ws.Cells.Item(89, 7).Value  = "Pepino dulce"
$ws.Cells.Item(89, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(89, 9).Value  = "Primera"
$ws.Cells.Item(89, 10).Value = 50
$ws.Cells.Item(89, 11).Value = 19000
$ws.Cells.Item(89, 12).Value = 19000
$ws.Cells.Item(89, 13).Value = 19000
$ws.Cells.Item(89, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(89, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(89, 16).Value = 1056
$ws.Cells.Item(89, 17).Value = 18
$ws.Cells.Item(89, 18).Value = "Hortaliza"
